$wb = $excel.ActiveWorkbook

# Worksheet used purely as a formatting template: it already has the
# "header row / data row" style pattern (s=9 header, s=2 body, s=5 for the
# 4th "Welcome Message" column) that the new test-data sheets reuse.
$styleSrc = $wb.Worksheets.Item("addImageURLtoItems")

# ---------------------------------------------------------------------
# New sheet 1: advSe033  (TC_ADV SEARCH_033)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws33 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws33.Name = "advSe033"

$styleSrc.Range("A1:E1").Copy()
$ws33.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:E2").Copy()
$ws33.Range("A2:E2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws33.Range("A1").Value = "TestCase ID"
$ws33.Range("B1").Value = "UserName"
$ws33.Range("C1").Value = "Password"
$ws33.Range("D1").Value = "Welcome Message"

$ws33.Range("A2").Value = "TC_ADV SEARCH_033"
$ws33.Range("B2").Value = "automationUser"
$ws33.Range("C2").Value = "unilog123##"
$ws33.Range("D2").Value = "Welcome, Automation !"
$ws33.Range("E2").Value = "Automation_PN"

$ws33.Range("E1").Value = "advSesearchinput"

$ws33.Range("A1:D2").Select()

# ---------------------------------------------------------------------
# New sheet 2: advSe039  (TC_ADV SEARCH_039)
# ---------------------------------------------------------------------
$ws34 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws33)
$ws34.Name = "advSe039"

$styleSrc.Range("A1:F1").Copy()
$ws34.Range("A1:F1").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:F2").Copy()
$ws34.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:F2").Copy()
$ws34.Range("A3:F3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws34.Range("A1").Value = "TestCase ID"
$ws34.Range("B1").Value = "UserName"
$ws34.Range("C1").Value = "Password"
$ws34.Range("D1").Value = "Welcome Message"
$ws34.Range("E1").Value = "Part Number"
$ws34.Range("F1").Value = "Item Name Template"

$ws34.Range("A2").Value = "TC_ADV SEARCH_039"
$ws34.Range("B2").Value = "automationUser"
$ws34.Range("C2").Value = "unilog123##"
$ws34.Range("D2").Value = "Welcome, Automation !"
$ws34.Range("E2").Value = "Automation_PN_1"
$ws34.Range("F2").Value = "Automation_PN_"

$ws34.Range("A3").Value = "TC_ADV SEARCH_039"
$ws34.Range("B3").Value = "automationUser"
$ws34.Range("C3").Value = "unilog123##"
$ws34.Range("D3").Value = "Welcome, Automation !"
$ws34.Range("E3").Value = "Automation_PN_2"
$ws34.Range("F3").Value = "Automation_PN_"

$ws34.Range("A1:F2").Select()

# ---------------------------------------------------------------------
# New sheet 3: advSearch_040  (TC_ADV SEARCH_040)
# ---------------------------------------------------------------------
$ws35 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws34)
$ws35.Name = "advSearch_040"

$styleSrc.Range("A1:F1").Copy()
$ws35.Range("A1:F1").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:F2").Copy()
$ws35.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:F2").Copy()
$ws35.Range("A3:F3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws35.Range("A1").Value = "TestCase ID"
$ws35.Range("B1").Value = "UserName"
$ws35.Range("C1").Value = "Password"
$ws35.Range("D1").Value = "Welcome Message"
$ws35.Range("E1").Value = "Part Number"
$ws35.Range("F1").Value = "Item Name Template"

$ws35.Range("A2").Value = "TC_ADV SEARCH_040"
$ws35.Range("B2").Value = "automationUser"
$ws35.Range("C2").Value = "unilog123##"
$ws35.Range("D2").Value = "Welcome, Automation !"
$ws35.Range("E2").Value = "Automation_PN_1"
$ws35.Range("F2").Value = "Automation_PN_"

$ws35.Range("A3").Value = "TC_ADV SEARCH_040"
$ws35.Range("B3").Value = "automationUser"
$ws35.Range("C3").Value = "unilog123##"
$ws35.Range("D3").Value = "Welcome, Automation !"
$ws35.Range("E3").Value = "Automation_PN_2"
$ws35.Range("F3").Value = "Automation_PN_"

$ws35.Range("F3").Select()
$excel.ActiveWindow.Zoom = 106

# ---------------------------------------------------------------------
# Tab-selection / view-state bookkeeping, matching the final state of
# the authored workbook: "addDocumenttoItems" is no longer the active
# sheet, its selection collapses to A1:D2 (active cell D2), and
# "addImageURLtoItems" becomes the active tab with the cursor left at F16.
# ---------------------------------------------------------------------
$wsDoc = $wb.Worksheets.Item("addDocumenttoItems")
$wsDoc.Activate()
$wsDoc.Range("A1:D2").Select()

$styleSrc.Activate()
$styleSrc.Range("F16").Select()
